$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("A4").Value = 0.6428759570744059
$ws.Range("B4").Value = 0.6444263053706816
$ws.Range("C4").Value = 0.6451317646088772
$ws.Range("D4").Value = 0.6406366763793544

# Row 6
$ws.Range("A6").Value = 0.7683220680197399
$ws.Range("D6").Value = 0.6403939541979988

# Row 7
$ws.Range("A7").Value = 0.7031906114146517
$ws.Range("D7").Value = 0.6401431656490685
